$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text values (e.g. "111.75", "42.483.72").
# Force a text number format on the whole data range before writing the new
# values so Excel does not silently reinterpret numeric-looking strings as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.483.72'
$ws.Range("E2").Value = '  -2.76%  '
$ws.Range("D3").Value = '2.226.20'
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").Value = '111.75'
$ws.Range("E5").Value = '  -7.43%  '
$ws.Range("D6").Value = '297.15'
$ws.Range("E6").Value = '  +11.11%  '
$ws.Range("E7").Value = '  -3.40%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  -3.85%  '
$ws.Range("D10").Value = '44.22'
$ws.Range("E10").Value = '  -8.77%  '
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  -3.60%  '
$ws.Range("D12").Value = '54.35'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '8.86'
$ws.Range("E13").Value = '  -4.66%  '
$ws.Range("D14").Value = '1.00'
$ws.Range("E14").Value = '  +8.57%  '
$ws.Range("D15").Value = '0.104'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").Value = '15.12'
$ws.Range("E16").Value = '  -3.13%  '
$ws.Range("D17").Value = '2.560.65'
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D18").Value = '2.233.83'
$ws.Range("E18").Value = '  -1.82%  '
$ws.Range("D19").Value = '42.491.99'
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").Value = '7.38'
$ws.Range("E20").Value = '  +5.92%  '
$ws.Range("E21").Value = '  -4.26%  '
$ws.Range("D22").Value = '72.91'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").Value = '3.52'
$ws.Range("E23").Value = '  +21.95%  '
$ws.Range("D24").Value = '2.36'
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("D25").Value = '229.68'
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").Value = '9.24'
$ws.Range("E26").Value = '  -4.36%  '
$ws.Range("D27").Value = '11.71'
$ws.Range("E27").Value = '  -3.01%  '
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").Value = '38.46'
$ws.Range("E30").Value = '  -10.04%  '
$ws.Range("D31").Value = '3.24'
$ws.Range("E31").Value = '  -4.27%  '
$ws.Range("D32").Value = '173.72'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").Value = '  -2.92%  '
$ws.Range("D34").Value = '0.0901'
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("D35").Value = '5.25'
$ws.Range("E35").Value = '  +13.64%  '
$ws.Range("E36").Value = '  -2.55%  '
$ws.Range("D37").Value = '4.33'
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("D39").Value = '0.0378'
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D41").Value = '2.43'
$ws.Range("E41").Value = '  -5.74%  '
$ws.Range("D42").Value = '72.23'
$ws.Range("E42").Value = '  -3.18%  '
$ws.Range("D43").Value = '0.236'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("E44").Value = '  -7.10%  '
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = '1.32'
$ws.Range("E46").Value = '  -5.70%  '
$ws.Range("E47").Value = '  -6.57%  '
$ws.Range("D48").Value = '1.34'
$ws.Range("E48").Value = '  +4.77%  '
$ws.Range("D49").Value = '103.39'
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '8.52'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '1.65'
$ws.Range("E51").Value = '  +6.39%  '
